# Konnect Bill Payment Verification Checks added
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value "Dango" into B6 (new shared string, referenced from row 6)
$ws.Range("B6").Value = "Dango"

# Move the active cell selection from A13 to A12
$ws.Range("A12").Select()
